# Update "Datos actualizados a 15 de Agosto de 2020 a las 18:20" -> "...19:37"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 19:37"

# --- Update country statistics with the latest reported figures ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5497613
$ws.Range("C4").Value = 21347
$ws.Range("D4").Value = 2878551
$ws.Range("E4").Value = 2447063
$ws.Range("G4").Value = 464
$ws.Range("H4").Value = 171999

# India (row 6)
$ws.Range("B6").Value = 2587461
$ws.Range("C6").Value = 62239
$ws.Range("D6").Value = 1858447
$ws.Range("E6").Value = 678934
$ws.Range("G6").Value = 946
$ws.Range("H6").Value = 50080

# Singapur (row 48)
$ws.Range("D48").Value = 51521
$ws.Range("E48").Value = 4113

# Marruecos (row 57)
$ws.Range("B57").Value = 41017
$ws.Range("C57").Value = 1776
$ws.Range("D57").Value = 28566
$ws.Range("E57").Value = 11819
$ws.Range("G57").Value = 21
$ws.Range("H57").Value = 632

# Suiza (row 58)
$ws.Range("D58").Value = 33200
$ws.Range("E58").Value = 2733

# Etiopia (row 67)
$ws.Range("B67").Value = 28894
$ws.Range("C67").Value = 1652
$ws.Range("D67").Value = 12037
$ws.Range("E67").Value = 16348
$ws.Range("G67").Value = 17
$ws.Range("H67").Value = 509

# Irlanda (row 68)
$ws.Range("B68").Value = 27191
$ws.Range("C68").Value = 196
$ws.Range("E68").Value = 2053

# Malaui (row 107)
$ws.Range("B107").Value = 5026
$ws.Range("C107").Value = 38
$ws.Range("D107").Value = 2623
$ws.Range("E107").Value = 2246
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 157

# Mozambique (row 125)
$ws.Range("B125").Value = 2791
$ws.Range("C125").Value = 83
$ws.Range("D125").Value = 1136
$ws.Range("E125").Value = 1636

# Tunez (row 136) - case count jumps ahead of Benin/Islandia/Sierra Leona
$ws.Range("B136").Value = 2023
$ws.Range("C136").Value = 120
$ws.Range("D136").Value = 1327
$ws.Range("E136").Value = 642
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 54

# Republica de Chipre (row 147) - case count jumps ahead of Georgia
$ws.Range("B147").Value = 1332
$ws.Range("C147").Value = 14
$ws.Range("D147").Value = 870
$ws.Range("E147").Value = 442
$ws.Range("H147").Value = 20

# Reunion (row 162)
$ws.Range("B162").Value = 816
$ws.Range("C162").Value = 40
$ws.Range("E162").Value = 154

# Mauricio (row 175)
$ws.Range("B175").Value = 346
$ws.Range("C175").Value = 1
$ws.Range("E175").Value = 2

# --- Re-sort the whole table by "Casos totales" (column B) descending, ---
# --- same as the live dashboard does every time the source data refreshes ---
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$dataRange.Sort($sortKey, 2, $null, $null, 2, $null, 1, 1)
